$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(24).EntireRow.Insert()
$ws.Rows.Item(24).EntireRow.Insert()

$ws.Range("A24").Value = 1
$ws.Range("B24").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C24").Value = 'Arica y Parinacota'
$ws.Range("D24").Value = 45012
$ws.Range("E24").Value = 15
$ws.Range("F24").Value = 100112028
$ws.Range("G24").Value = 'Sandia'
$ws.Range("H24").Value = 'Sin especificar'
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 500
$ws.Range("K24").Value = 450
$ws.Range("L24").Value = 480
$ws.Range("M24").Value = 471
$ws.Range("N24").Value = '$/kilo (volumen en unidades)'
$ws.Range("O24").Value = 'Perú'
$ws.Range("P24").Value = 471
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = 'Hortaliza'

$ws.Range("A25").Value = 1
$ws.Range("B25").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C25").Value = 'Arica y Parinacota'
$ws.Range("D25").Value = 45012
$ws.Range("E25").Value = 15
$ws.Range("F25").Value = 100112028
$ws.Range("G25").Value = 'Sandia'
$ws.Range("H25").Value = 'Sin especificar'
$ws.Range("I25").Value = 'Segunda'
$ws.Range("J25").Value = 300
$ws.Range("K25").Value = 450
$ws.Range("L25").Value = 480
$ws.Range("M25").Value = 465
$ws.Range("N25").Value = '$/kilo (volumen en unidades)'
$ws.Range("O25").Value = 'Perú'
$ws.Range("P25").Value = 465
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = 'Hortaliza'
